# feat(EM): add whitelist & important only mode
# Adds three new localization rows (ids + JP/EN text) to the "General"
# sheet, right after the existing "AI Service Status" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 113: em_ui_global_cooldown
# (C/D share the same literal string, so write order doesn't matter here)
$ws.Range("A113").Value = "em_ui_global_cooldown"
$ws.Range("C113").Value = "Global Request Cooldown(s)"
$ws.Range("D113").Value = "Global Request Cooldown(s)"

# Row 114: em_ui_tab_whitelist
# Write D (EN) before C (JP) so new shared-string indices come out in the
# same order as the source workbook.
$ws.Range("A114").Value = "em_ui_tab_whitelist"
$ws.Range("D114").Value = "Whitelist"
$ws.Range("C114").Value = "ホワイトリスト"

# Row 115: em_ui_whitelist
$ws.Range("A115").Value = "em_ui_whitelist"
$ws.Range("D115").Value = "AI Service Whitelist Mode: {0}"
$ws.Range("C115").Value = "AIサービスホワイトリストモード：{0}"

# Keep the active selection consistent with the edited area (mirrors the
# author scrolling to / selecting the newly-added rows after editing).
$ws.Range("C117").Select()
